$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the first debitor sample row (row 2) with the new draft data:
# debitor_ident (A), Sak_Nr (B), FNR (C)
$ws.Range("A2").Value = "13088334935"
$ws.Range("B2").Value = "267794"
$ws.Range("C2").Value = "13088334935"

# The remaining leftover sample rows (3-6) are no longer needed for this draft
$ws.Range("A3:C6").Value = ""

# Leave the selection on the data that was just filled in
$ws.Range("A2:C3").Select()
